$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the hidden "_GoBack" bookmark from paragraph 1.
#    ("_GoBack" is a hidden bookmark, so it is excluded from
#    Bookmarks.Count / enumeration, but Exists/Item still resolve
#    it by name.)
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------
# 2. Merge the split runs in paragraph 3 ("Add some comments...")
#    into a single run with the full updated sentence.
# ---------------------------------------------------------------
$oldP3 = "Add some comments about Version management after this line, " + `
    "or just add some text so there is a change to this file.  " + `
    "Remember that your GitHub user " + "ID" + " must be submitted in you" + "r" + " assignment report!"
$newP3 = "Add some comments about Version management after this line, " + `
    "or just add some text so there is a change to this file.  " + `
    "Remember that your GitHub user ID must be submitted in your assignment report!"
$d.Content.Find.Execute($oldP3, $true, $false, $false, $false, $false, $true, 1, $false, $newP3, 2)

# ---------------------------------------------------------------
# 3. Merge the split runs in paragraph 4 (">>> your stuff...")
#    into a single run, dropping the proofErr markers.
# ---------------------------------------------------------------
$oldP4 = ">>" + ">  your" + " stuff after this line >>>"
$newP4 = ">>> your stuff after this line >>>"
$d.Content.Find.Execute($oldP4, $true, $false, $false, $false, $false, $true, 1, $false, $newP4, 2)

# ---------------------------------------------------------------
# 4. Replace "Ben changing things up!" with the new closing text.
# ---------------------------------------------------------------
$d.Content.Find.Execute("Ben changing things up!", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Configuration management system is a system that process and maintaining performance product for functionality and design by user", 2)

# ---------------------------------------------------------------
# 5. Drop the trailing empty paragraphs at the end of the body.
#    Word will not allow deleting the document's very last
#    paragraph mark, so instead delete the mark that belongs to
#    the paragraph just *before* the last one -- this merges the
#    last (blank) paragraph upward and shrinks the paragraph
#    count by one. Repeat while the document still ends in a
#    blank paragraph and there is a preceding paragraph to
#    absorb it into.
# ---------------------------------------------------------------
while ($d.Paragraphs.Count -gt 1) {
    $count = $d.Paragraphs.Count
    $last = $d.Paragraphs.Item($count)
    if ($last.Range.Text.Trim().Length -eq 0) {
        $secondLast = $d.Paragraphs.Item($count - 1)
        $r = $d.Range($secondLast.Range.End - 1, $secondLast.Range.End)
        $r.Delete()
    } else {
        break
    }
}

# ---------------------------------------------------------------
# 6. Justify (both) the first four paragraphs.
# ---------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $d.Paragraphs.Item($i).Alignment = 3
}

Write-Output "done"
